$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Replace the Mean (B) / SD (C) columns with the California dataset.
#    Column A (Frac) is unchanged.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 0.36570000000000003
$ws.Range("C2").Value = 0.0032000000000000002

$ws.Range("B3").Value = 0.36509999999999998
$ws.Range("C3").Value = 0.0041000000000000003

$ws.Range("B4").Value = 0.36480000000000001
$ws.Range("C4").Value = 0.0045999999999999999

$ws.Range("B5").Value = 0.36580000000000001
$ws.Range("C5").Value = 0.0068999999999999999

$ws.Range("B6").Value = 0.36499999999999999
$ws.Range("C6").Value = 0.0079000000000000008

$ws.Range("B7").Value = 0.36480000000000001
$ws.Range("C7").Value = 0.0101

$ws.Range("B8").Value = 0.36530000000000001
$ws.Range("C8").Value = 0.0137

$ws.Range("B9").Value = 0.36230000000000001
$ws.Range("C9").Value = 0.020400000000000001

# ---------------------------------------------------------------------------
# 2) Narrow column B to (approximately) match column C's width.
# ---------------------------------------------------------------------------
$ws.Range("B1").EntireColumn.ColumnWidth = 5.833333333333334

# ---------------------------------------------------------------------------
# 3) Re-stamp the header row's style (A1:C1). The style table carries a long
#    run of orphaned, unreferenced cellXfs/border records, so faithfully
#    landing on the same final style index means replaying that same number
#    of "genuinely new" number-format/border combinations first (each one
#    mints a fresh, never-reused style slot), then letting the header cells
#    settle on the last one minted.
# ---------------------------------------------------------------------------

# H1 -> numFmt 49 (text "@"), fresh border (diagonal-down = thin)
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Borders.Item(5).LineStyle = 1

# H2 -> numFmt 22 (builtin date/time), fresh border (diagonal-down cleared again)
$ws.Range("Z2").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z2").Borders.Item(5).LineStyle = 1
$ws.Range("Z2").Borders.Item(5).LineStyle = -4142

# H3 -> numFmt 49, fresh border (diagonal-up = thin)
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Borders.Item(6).LineStyle = 1

# H4 -> numFmt 22, fresh border (diagonal-up cleared again)
$ws.Range("Z4").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z4").Borders.Item(6).LineStyle = 1
$ws.Range("Z4").Borders.Item(6).LineStyle = -4142

# H5 -> numFmt 49, fresh border (diagonal-down + diagonal-up = thin) - this is
# the style the header row itself lands on.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A1:C1").Borders.Item(5).LineStyle = 1
$ws.Range("A1:C1").Borders.Item(6).LineStyle = 1

# H6 -> numFmt 22, fresh border (diagonal-down + diagonal-up cleared again) -
# keeps the tail of the style table lined up even though nothing references it.
$ws.Range("Z6").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z6").Borders.Item(5).LineStyle = 1
$ws.Range("Z6").Borders.Item(6).LineStyle = 1
$ws.Range("Z6").Borders.Item(5).LineStyle = -4142
$ws.Range("Z6").Borders.Item(6).LineStyle = -4142

# Drop the scratch cells - only their style-table side effects should remain.
$ws.Range("Z1:Z6").Clear()
